$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lama2"
$ws.Range("C2").Value = "Dag1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.608931666666666
$ws.Range("H2").Value = 13.826795
$ws.Range("I2").Value = 0.02269509467890621
$ws.Range("J2").Value = 0.02269509467890622
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 18.76192366666666
$ws.Range("N2").Value = 56.285771
$ws.Range("O2").Value = 0.1222461152048115
$ws.Range("P2").Value = 0.1222461152048115
$ws.Range("Q2").Value = 86.47242411488276
$ws.Range("R2").Value = 778.2518170339449
$ws.Range("S2").Value = 0.002774387158701674
$ws.Range("T2").Value = 0.002774387158701675

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lama2"
$ws.Range("C3").Value = "Dag1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.608931666666666
$ws.Range("H3").Value = 13.826795
$ws.Range("I3").Value = 0.02269509467890621
$ws.Range("J3").Value = 0.02269509467890622
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 38.59812166666666
$ws.Range("N3").Value = 115.794365
$ws.Range("O3").Value = 0.2514918252404857
$ws.Range("P3").Value = 0.2514918252404857
$ws.Range("Q3").Value = 177.8961052233527
$ws.Range("R3").Value = 1601.064947010175
$ws.Range("S3").Value = 0.005707630784803758
$ws.Range("T3").Value = 0.005707630784803759

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lama2"
$ws.Range("C4").Value = "Dag1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.608931666666666
$ws.Range("H4").Value = 13.826795
$ws.Range("I4").Value = 0.02269509467890621
$ws.Range("J4").Value = 0.02269509467890622
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.691683333333334
$ws.Range("N4").Value = 20.07505
$ws.Range("O4").Value = 0.04360066196912097
$ws.Range("P4").Value = 0.04360066196912097
$ws.Range("Q4").Value = 30.84151121830555
$ws.Range("R4").Value = 277.57360096475
$ws.Range("S4").Value = 0.0009895211514521858
$ws.Range("T4").Value = 0.000989521151452186

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lama2"
$ws.Range("C5").Value = "Dag1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.608931666666666
$ws.Range("H5").Value = 13.826795
$ws.Range("I5").Value = 0.02269509467890621
$ws.Range("J5").Value = 0.02269509467890622
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 89.42491666666666
$ws.Range("N5").Value = 268.27475
$ws.Range("O5").Value = 0.5826613975855819
$ws.Range("P5").Value = 0.5826613975855818
$ws.Range("Q5").Value = 412.1533302140277
$ws.Range("R5").Value = 3709.379971926249
$ws.Range("S5").Value = 0.0132235555839486
$ws.Range("T5").Value = 0.0132235555839486

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lama2"
$ws.Range("C6").Value = "Dag1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 168.218394
$ws.Range("H6").Value = 504.655182
$ws.Range("I6").Value = 0.8283334739316415
$ws.Range("J6").Value = 0.8283334739316416
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.76192366666666
$ws.Range("N6").Value = 56.285771
$ws.Range("O6").Value = 0.1222461152048115
$ws.Range("P6").Value = 0.1222461152048115
$ws.Range("Q6").Value = 3156.100667557258
$ws.Range("R6").Value = 28404.90600801532
$ws.Range("S6").Value = 0.1012605492822492
$ws.Range("T6").Value = 0.1012605492822492

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lama2"
$ws.Range("C7").Value = "Dag1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 168.218394
$ws.Range("H7").Value = 504.655182
$ws.Range("I7").Value = 0.8283334739316415
$ws.Range("J7").Value = 0.8283334739316416
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 38.59812166666666
$ws.Range("N7").Value = 115.794365
$ws.Range("O7").Value = 0.2514918252404857
$ws.Range("P7").Value = 0.2514918252404857
$ws.Range("Q7").Value = 6492.91403818327
$ws.Range("R7").Value = 58436.22634364943
$ws.Range("S7").Value = 0.2083190972668608
$ws.Range("T7").Value = 0.2083190972668608

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Lama2"
$ws.Range("C8").Value = "Dag1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 168.218394
$ws.Range("H8").Value = 504.655182
$ws.Range("I8").Value = 0.8283334739316415
$ws.Range("J8").Value = 0.8283334739316416
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.691683333333334
$ws.Range("N8").Value = 20.07505
$ws.Range("O8").Value = 0.04360066196912097
$ws.Range("P8").Value = 0.04360066196912097
$ws.Range("Q8").Value = 1125.6642234899
$ws.Range("R8").Value = 10130.9780114091
$ws.Range("S8").Value = 0.03611588779460118
$ws.Range("T8").Value = 0.03611588779460118

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Lama2"
$ws.Range("C9").Value = "Dag1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 168.218394
$ws.Range("H9").Value = 504.655182
$ws.Range("I9").Value = 0.8283334739316415
$ws.Range("J9").Value = 0.8283334739316416
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 89.42491666666666
$ws.Range("N9").Value = 268.27475
$ws.Range("O9").Value = 0.5826613975855819
$ws.Range("P9").Value = 0.5826613975855818
$ws.Range("Q9").Value = 15042.9158652505
$ws.Range("R9").Value = 135386.2427872545
$ws.Range("S9").Value = 0.4826379395879304
$ws.Range("T9").Value = 0.4826379395879303

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Lama2"
$ws.Range("C10").Value = "Dag1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1627236666666667
$ws.Range("H10").Value = 0.488171
$ws.Range("I10").Value = 0.00080127658394417
$ws.Range("J10").Value = 0.00080127658394417
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 18.76192366666666
$ws.Range("N10").Value = 56.285771
$ws.Range("O10").Value = 0.1222461152048115
$ws.Range("P10").Value = 0.1222461152048115
$ws.Range("Q10").Value = 3.053009012760111
$ws.Range("R10").Value = 27.477081114841
$ws.Range("S10").Value = 0.00009795294959175684
$ws.Range("T10").Value = 0.00009795294959175684

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Lama2"
$ws.Range("C11").Value = "Dag1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1627236666666667
$ws.Range("H11").Value = 0.488171
$ws.Range("I11").Value = 0.00080127658394417
$ws.Range("J11").Value = 0.00080127658394417
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 38.59812166666666
$ws.Range("N11").Value = 115.794365
$ws.Range("O11").Value = 0.2514918252404857
$ws.Range("P11").Value = 0.2514918252404857
$ws.Range("Q11").Value = 6.280827884046111
$ws.Range("R11").Value = 56.527450956415
$ws.Range("S11").Value = 0.0002015145106185806
$ws.Range("T11").Value = 0.0002015145106185806

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Lama2"
$ws.Range("C12").Value = "Dag1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1627236666666667
$ws.Range("H12").Value = 0.488171
$ws.Range("I12").Value = 0.00080127658394417
$ws.Range("J12").Value = 0.00080127658394417
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 6.691683333333334
$ws.Range("N12").Value = 20.07505
$ws.Range("O12").Value = 0.04360066196912097
$ws.Range("P12").Value = 0.04360066196912097
$ws.Range("Q12").Value = 1.088895248172222
$ws.Range("R12").Value = 9.800057233550001
$ws.Range("S12").Value = 0.00003493618948032174
$ws.Range("T12").Value = 0.00003493618948032174

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Lama2"
$ws.Range("C13").Value = "Dag1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1627236666666667
$ws.Range("H13").Value = 0.488171
$ws.Range("I13").Value = 0.00080127658394417
$ws.Range("J13").Value = 0.00080127658394417
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 89.42491666666666
$ws.Range("N13").Value = 268.27475
$ws.Range("O13").Value = 0.5826613975855819
$ws.Range("P13").Value = 0.5826613975855818
$ws.Range("Q13").Value = 14.55155033136111
$ws.Range("R13").Value = 130.96395298225
$ws.Range("S13").Value = 0.0004668729342535109
$ws.Range("T13").Value = 0.0004668729342535108

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Lama2"
$ws.Range("C14").Value = "Dag1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 30.09047233333333
$ws.Range("H14").Value = 90.271417
$ws.Range("I14").Value = 0.148170154805508
$ws.Range("J14").Value = 0.148170154805508
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 18.76192366666666
$ws.Range("N14").Value = 56.285771
$ws.Range("O14").Value = 0.1222461152048115
$ws.Range("P14").Value = 0.1222461152048115
$ws.Range("Q14").Value = 564.5551450119452
$ws.Range("R14").Value = 5080.996305107506
$ws.Range("S14").Value = 0.0181132258142689
$ws.Range("T14").Value = 0.0181132258142689

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Lama2"
$ws.Range("C15").Value = "Dag1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 30.09047233333333
$ws.Range("H15").Value = 90.271417
$ws.Range("I15").Value = 0.148170154805508
$ws.Range("J15").Value = 0.148170154805508
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 38.59812166666666
$ws.Range("N15").Value = 115.794365
$ws.Range("O15").Value = 0.2514918252404857
$ws.Range("P15").Value = 0.2514918252404857
$ws.Range("Q15").Value = 1161.435712129467
$ws.Range("R15").Value = 10452.9214091652
$ws.Range("S15").Value = 0.03726358267820254
$ws.Range("T15").Value = 0.03726358267820254

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Lama2"
$ws.Range("C16").Value = "Dag1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 30.09047233333333
$ws.Range("H16").Value = 90.271417
$ws.Range("I16").Value = 0.148170154805508
$ws.Range("J16").Value = 0.148170154805508
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 6.691683333333334
$ws.Range("N16").Value = 20.07505
$ws.Range("O16").Value = 0.04360066196912097
$ws.Range("P16").Value = 0.04360066196912097
$ws.Range("Q16").Value = 201.3559122050945
$ws.Range("R16").Value = 1812.20320984585
$ws.Range("S16").Value = 0.006460316833587281
$ws.Range("T16").Value = 0.006460316833587281

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Lama2"
$ws.Range("C17").Value = "Dag1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 30.09047233333333
$ws.Range("H17").Value = 90.271417
$ws.Range("I17").Value = 0.148170154805508
$ws.Range("J17").Value = 0.148170154805508
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 89.42491666666666
$ws.Range("N17").Value = 268.27475
$ws.Range("O17").Value = 0.5826613975855819
$ws.Range("P17").Value = 0.5826613975855818
$ws.Range("Q17").Value = 2690.837980868972
$ws.Range("R17").Value = 24217.54182782075
$ws.Range("S17").Value = 0.08633302947944935
$ws.Range("T17").Value = 0.08633302947944932

